# Update "想去人数" (number of people interested) counts.
# F2: 347 -> 348, F5: 290 -> 291, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 348
    $ws.Range("F5").Value = 291
}
